$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.607.79'
$ws.Range('E2').Value = '  -7.34%  '
$ws.Range('D3').Value = '1.696.50'
$ws.Range('E3').Value = '  -6.05%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.70'
$ws.Range('E5').Value = '  -5.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5141'
$ws.Range('E6').Value = '  -13.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.006'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2650'
$ws.Range('E8').Value = '  -4.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '22.23'
$ws.Range('E9').Value = '  -4.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06284'
$ws.Range('E10').Value = '  -7.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07355'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').Value = '1.703.95'
$ws.Range('E12').Value = '  -5.69%  '
$ws.Range('E13').Value = '  -4.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5856'
$ws.Range('E14').Value = '  -6.04%  '
$ws.Range('D15').Value = '1.927.15'
$ws.Range('E15').Value = '  -6.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008406'
$ws.Range('E16').Value = '  -8.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.57'
$ws.Range('E17').Value = '  -13.35%  '
$ws.Range('D18').Value = '26.650.39'
$ws.Range('E18').Value = '  -7.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.031'
$ws.Range('E19').Value = '  -8.13%  '
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.95'
$ws.Range('E21').Value = '  -4.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '186.98'
$ws.Range('E22').Value = '  -11.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.283'
$ws.Range('E23').Value = '  -8.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.007'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.00'
$ws.Range('E25').Value = '  -5.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.585'
$ws.Range('E26').Value = '  -3.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1153'
$ws.Range('E27').Value = '  -8.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.77'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.325'
$ws.Range('E29').Value = '  -7.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05698'
$ws.Range('E30').Value = '  -7.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.341'
$ws.Range('E31').Value = '  -6.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.520'
$ws.Range('E32').Value = '  -7.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.513'
$ws.Range('E33').Value = '  -6.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.655'
$ws.Range('E34').Value = '  -4.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.026'
$ws.Range('E35').Value = '  -3.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6036'
$ws.Range('E36').Value = '  -6.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.373'
$ws.Range('E37').Value = '  -5.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.685'
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('D39').Value = '1.103.69'
$ws.Range('E39').Value = '  -3.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01604'
$ws.Range('E40').Value = '  -5.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8636'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.865'
$ws.Range('E42').Value = '  -10.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.005'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.04'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('D45').Value = '1.856.24'
$ws.Range('E45').Value = '  -5.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000112'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.76'
$ws.Range('E47').Value = '  -6.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.199'
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.002'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05241'
$ws.Range('E50').Value = '  -4.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4328'
$ws.Range('E51').Value = '  -3.37%  '
